$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report generated timestamp
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:01 AM"

# Update total billed amount
$ws.Range("C8").Value = 216.17

# Clear the Scope ID # value (was "#NO MATCH", now empty)
$ws.Range("G10").Value = ""

# Update pricing for line item and total
$ws.Range("H16").Value = 216.17
$ws.Range("H17").Value = 216.17
